$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column widths for J, K, L (closest achievable to 17.75 / 22.25 / 25.625) ---
$ws.Range("J1").ColumnWidth = 17
$ws.Range("K1").ColumnWidth = 21.5
$ws.Range("L1").ColumnWidth = 24.9

# --- Extra data point appended to the averages row ---
$ws.Range("G18").Value = 4.0016499999999997

# --- Header row for the Taylor-series expansion table ---
$ws.Range("J22").Value = "Term 0 (1 + x)"
$ws.Range("K22").Value = "Term 1 ( + x^2 * 0.5)"
$ws.Range("L22").Value = "Term 2 (+x^3 / 6)"

# --- Styled (but empty) label cell above the table ---
$ws.Range("D23").RowHeight = 22.5
$ws.Range("D23").Font.Name = "JetBrains Mono"
$ws.Range("D23").Font.Size = 17.3
$ws.Range("D23").Font.Color = 16758905
$ws.Range("D23").Font.Family = 3

# --- x values used by the Taylor series approximation of exp(x) ---
$ws.Range("I23").Value = 0.4653
$ws.Range("I24").Value = 0.7761
$ws.Range("I25").Value = 0.7094
$ws.Range("I26").Value = 0.5759
$ws.Range("I27").Value = 0.1992
$ws.Range("I28").Value = 0.788
$ws.Range("I29").Value = 0.4316
$ws.Range("I30").Value = 0.6491

# --- Row 23: plain (non-shared) formulas ---
$ws.Range("J23").Formula = "=1+I23"
$ws.Range("K23").Formula = "=J23  + I23 * I23 * 0.5"
$ws.Range("L23").Formula = "=K23 +POWER(I23, 3)"
$ws.Range("M23").Formula = "=POWER(I23, 3)"

# --- Rows 24:30: shared formulas ---
$ws.Range("J24:J30").Formula = "=1+I24"
$ws.Range("K24:K30").Formula = "=J24  + I24 * I24 * 0.5"
$ws.Range("L24:L30").Formula = "=K24 +POWER(I24, 3) / 6"
$ws.Range("M24:M30").Formula = "=POWER(I24, 3)"

# --- View state: scroll the sheet and move the active selection ---
$ws.Range("N25").Select()
